# order list원본과 sales list 원본을 넣어둠
#
# 1) Swap the "Data" and "Clients" sheet names (the physical sheet behind
#    the old "Data" tab becomes "Clients", and vice versa - sheetId/r:id
#    stay attached to the same physical sheet, only the tab names trade
#    places).
# 2) Replace the (now-"Clients") sheet content with the client record for
#    콕스 (16 columns, 2 rows).
# 3) Replace the (now-"Data") sheet content with two quote line items for
#    콕스 (30 columns, 3 rows).
# 4) Append two rows to the Log sheet describing the above actions.

$wb = $excel.ActiveWorkbook

function Set-CellValue($sheet, $row, $col, $value) {
    if ($null -eq $value) {
        return
    }
    if ($value -is [string] -and $value -match '^\d{4}-\d{2}-\d{2}$') {
        # Plain yyyy-MM-dd strings must stay literal text, not silently
        # turn into an Excel date serial number.
        $sheet.Cells.Item($row, $col).NumberFormat = "@"
    }
    $sheet.Cells.Item($row, $col).Value = $value
}

function Set-RowValues($sheet, $row, $values) {
    $col = 1
    foreach ($v in $values) {
        Set-CellValue $sheet $row $col $v
        $col = $col + 1
    }
}

# --- Step 1: swap sheet names --------------------------------------------
$wsOldData = $wb.Worksheets.Item("Data")
$wsOldClients = $wb.Worksheets.Item("Clients")

$wsOldData.Name = "__TMP__"
$wsOldClients.Name = "Data"
$wsOldData.Name = "Clients"

$wsClients = $wb.Worksheets.Item("Clients")
$wsData = $wb.Worksheets.Item("Data")
$wsLog = $wb.Worksheets.Item("Log")
$wsFormatSource = $wb.Worksheets.Item("Payment")

# --- Step 2: Clients sheet content ----------------------------------------
# The physical sheet previously held the 26-column "Data" header, so wipe
# it completely before laying out the new, narrower Clients table.
$wsClients.Cells.Clear()

$clientsHeader = @("업체명","사업자번호","대표자","전화번호","이메일","주소","특이사항","운송방법","운송계정","국가","통화","담당자","수출허가구분","수출허가번호","수출허가만료일","사업자등록증경로")
Set-RowValues $wsClients 1 $clientsHeader

$clientsRow2 = @("콕스","-","-","010-2314-1234","sue@coxcamera.com","디지털로 242","ㅁㅇㄴㄻㄴㅇㄹ","DHL","콕","KR","KRW","하수민","해당 없음","해당 없음","2025-12-31","//cox_biz/business/SalesManager/attachments\사업자등록증\사업자등록증_콕스__251209.png")
Set-RowValues $wsClients 2 $clientsRow2

# Reapply the bold/centered/bordered header style (same one used on every
# other sheet's row 1) to the newly written header row.
$wsFormatSource.Range("A1").Copy()
$wsClients.Range("A1:P1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Step 3: Data sheet content --------------------------------------------
# The physical sheet previously held the 9-column "Clients" header, so wipe
# it completely before laying out the new, wider Data table.
$wsData.Cells.Clear()

$dataHeader = @("관리번호","업체명","모델명","Description","수량","단가","환율","세율(%)","공급가액","세액","합계금액","기수금액","미수금액","견적일","수주일","출고예정일","출고일","선적일","입금완료일","세금계산서발행일","계산서번호","수출신고번호","수출신고필증경로","Status","비고","주문요청사항","구분","프로젝트명","통화","품목명")
Set-RowValues $wsData 1 $dataHeader

$dataRow2 = @("QT-251209-001","콕스","ㅁㄴㅇㄹ","ㄴㅇㄻㄴㅇㄹ",1,1000,1,10,1000,100,1100,0,1100,"2025-12-09",$null,$null,$null,$null,$null,$null,$null,$null,$null,"견적","ㄻㄴㅇㄻㄹㅁㄴ",$null,"내수","ㅁㄴㅇㄻㄴㅇ","KRW","ㅇㄹㅁㄴㅇㄹ")
Set-RowValues $wsData 2 $dataRow2

$dataRow3 = @("QT-251209-001","콕스","ㅇㄻㄴ","ㅇㄻㄴㅇㄹ",1,1200,1,10,1200,120,1320,0,1320,"2025-12-09",$null,$null,$null,$null,$null,$null,$null,$null,$null,"견적","ㄻㄴㅇㄻㄹㅁㄴ",$null,"내수","ㅁㄴㅇㄻㄴㅇ","KRW","ㅁㄴㅇㄻ")
Set-RowValues $wsData 3 $dataRow3

$wsFormatSource.Range("A1").Copy()
$wsData.Range("A1:AD1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Step 4: Log sheet content ----------------------------------------------
$logRow2 = @("2025-12-09 15:57:24","sue","업체 등록","업체명: 콕스")
Set-RowValues $wsLog 2 $logRow2

$logRow3 = @("2025-12-09 15:57:46","sue","견적 등록","번호 [QT-251209-001] / 업체 [콕스]")
Set-RowValues $wsLog 3 $logRow3
